$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price / link / volume(1h) data.
# Numeric-looking price strings must be forced to Text format
# before assignment, otherwise Excel auto-converts them to
# numbers (stripping the thousands-separator dots etc.).

# Row 2: Bitcoin
$ws.Range("D2").Value = "41.481.18"
$ws.Range("E2").Value = "  -2.28%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.457.28"

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.88%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.66"
$ws.Range("E5").Value = "  -0.95%  "

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "90.35"
$ws.Range("E6").Value = "  -7.24%  "

# Row 7: XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.536"
$ws.Range("E7").Value = "  -4.50%  "

# Row 8: USDC
$ws.Range("E8").Value = "  +0.76%  "

# Row 9: Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.484"
$ws.Range("E9").Value = "  -6.82%  "

# Row 10: Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.45"
$ws.Range("E10").Value = "  -7.33%  "

# Row 11: Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0770"
$ws.Range("E11").Value = "  -3.58%  "

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.837.68"
$ws.Range("E13").Value = "  -2.32%  "

# Row 14: Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.75"
$ws.Range("E14").Value = "  -6.08%  "

# Row 15: WrappedEther
$ws.Range("D15").Value = "2.550.00"
$ws.Range("E15").Value = "  +0.58%  "

# Row 16: Chainlink
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.06"
$ws.Range("E16").Value = "  +0.01%  "

# Row 17: Polygon
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.767"
$ws.Range("E17").Value = "  -5.47%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "41.192.55"
$ws.Range("E18").Value = "  -3.04%  "

# Row 19: Uniswap
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.20"
$ws.Range("E19").Value = "  -5.29%  "

# Row 20: ShibaInu
$ws.Range("D20").Value = "0.0₃0907"
$ws.Range("E20").Value = "  -3.35%  "

# Row 21: Litecoin
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.90"
$ws.Range("E21").Value = "  +1.00%  "

# Row 22: InternetComputer(DFINITY)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.84"
$ws.Range("E22").Value = "  -9.91%  "

# Row 23: BitcoinCash
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.06"
$ws.Range("E23").Value = "  -3.00%  "

# Row 24: PancakeSwap
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.69"
$ws.Range("E24").Value = "  -5.56%  "

# Row 25: Dai
$ws.Range("E25").Value = "  -0.01%  "

# Row 26: ImmutableX
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.85"
$ws.Range("E26").Value = "  -6.58%  "

# Row 27: EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.76"
$ws.Range("E27").Value = "  -6.27%  "

# Row 28: Toncoin
$ws.Range("E28").Value = "  -0.37%  "

# Row 29: Cosmos
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.57"
$ws.Range("E29").Value = "  -3.73%  "

# Row 30: InjectiveProtocol
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.30"
$ws.Range("E30").Value = "  -5.70%  "

# Row 31: Monero
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "151.86"
$ws.Range("E31").Value = "  -2.34%  "

# Row 32: Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.33"
$ws.Range("E32").Value = "  -8.48%  "

# Row 33: ApeXProtocol
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.57"
$ws.Range("E33").Value = "  -4.42%  "

# Row 34: WEMIXToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.55"
$ws.Range("E34").Value = "  -3.33%  "

# Row 35: Hedera
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0743"
$ws.Range("E35").Value = "  -4.60%  "

# Row 36: Celestia
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.14"
$ws.Range("E36").Value = "  -1.02%  "

# Row 37: LidoDAOToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.94"
$ws.Range("E37").Value = "  -5.39%  "

# Row 38: ARBITRUM
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.81"
$ws.Range("E38").Value = "  -7.99%  "

# Row 39: Stellar
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.111"
$ws.Range("E39").Value = "  -4.08%  "

# Row 40: Kaspa
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0985"
$ws.Range("E40").Value = "  -8.70%  "

# Row 41: RenderToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.96"
$ws.Range("E41").Value = "  -6.23%  "

# Row 42: FirstDigitalUSD
$ws.Range("E42").Value = "  +1.11%  "

# Row 43: EnergySwap
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.21"
$ws.Range("E43").Value = "  -8.76%  "

# Row 44: Maker
$ws.Range("D44").Value = "1.945.71"
$ws.Range("E44").Value = "  -3.98%  "

# Row 45: VeChain
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0277"
$ws.Range("E45").Value = "  -5.57%  "

# Row 46: NEARProtocol
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.89"
$ws.Range("E46").Value = "  -9.21%  "

# Row 47: FraxShare -> RocketPoolETH (rows 47/48 swapped)
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "2.712.59"
$ws.Range("E47").Value = "  -1.86%  "

# Row 48: RocketPoolETH -> FraxShare (rows 47/48 swapped)
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.55"
$ws.Range("E48").Value = "  -3.11%  "

# Row 49: Aave
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "94.84"
$ws.Range("E49").Value = "  -4.87%  "

# Row 50: ordi
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.77"
$ws.Range("E50").Value = "  -6.88%  "

# Row 51: Algorand
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.173"
$ws.Range("E51").Value = "  -7.17%  "

